$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Records screened (n1 = 186)" -> "... = 187)" ---
$sh1 = $s.Shapes.Item(3)
$tr1 = $sh1.TextFrame.TextRange
$run186 = $tr1.Find("186")
$run186.Text = "187"

# --- "Records excluded because ... (n1 = 173)" -> "... = 174)" ---
$sh2 = $s.Shapes.Item(4)
$tr2 = $sh2.TextFrame.TextRange
$run173 = $tr2.Find("= 173)")
$run173.Text = "= 174)"
